# Re-process the sheet with the newly curated dimensions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: concept/dimension references
$ws.Range("A2").Value = "iaest-measure:tipo-de-hogar-2"
$ws.Range("C2").Value = "iaest-measure:tipo-hogar-1"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"

# Row 3: role (dim / medida)
$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "dim"

# Row 4: datatype / URI role
$ws.Range("A4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("F4").Value = "URI-Comunidad"

# Row 5: no longer used - remove the mapping file references entirely
$ws.Range("A5:H5").EntireRow.Delete()
